{"js": "// Apply the text replacements described by the diff: the header date and\n// the 25 three-digit-by-one-digit multiplication equations in the table.\nconst replacements = [\n  [\"2025-03-15 Saturday\", \"2025-03-16 Sunday\"],\n  [\"229\u00d72=458\", \"457\u00d74=1828\"],\n  [\"981\u00d77=6867\", \"988\u00d75=4940\"],\n  [\"730\u00d76=4380\", \"393\u00d77=2751\"],\n  [\"588\u00d79=5292\", \"157\u00d74=628\"],\n  [\"324\u00d72=648\", \"357\u00d72=714\"],\n  [\"815\u00d75=4075\", \"143\u00d78=1144\"],\n  [\"920\u00d72=1840\", \"916\u00d75=4580\"],\n  [\"867\u00d72=1734\", \"461\u00d74=1844\"],\n  [\"909\u00d72=1818\", \"458\u00d75=2290\"],\n  [\"765\u00d77=5355\", \"616\u00d75=3080\"],\n  [\"654\u00d72=1308\", \"859\u00d76=5154\"],\n  [\"515\u00d73=1545\", \"329\u00d72=658\"],\n  [\"238\u00d77=1666\", \"425\u00d72=850\"],\n  [\"450\u00d72=900\", \"216\u00d78=1728\"],\n  [\"805\u00d73=2415\", \"245\u00d73=735\"],\n  [\"277\u00d75=1385\", \"276\u00d77=1932\"],\n  [\"291\u00d73=873\", \"344\u00d76=2064\"],\n  [\"651\u00d73=1953\", \"136\u00d77=952\"],\n  [\"116\u00d77=812\", \"516\u00d78=4128\"],\n  [\"472\u00d72=944\", \"337\u00d73=1011\"],\n  [\"628\u00d76=3768\", \"188\u00d75=940\"],\n  [\"496\u00d79=4464\", \"136\u00d74=544\"],\n  [\"227\u00d77=1589\", \"972\u00d78=7776\"],\n  [\"582\u00d76=3492\", \"920\u00d77=6440\"],\n  [\"421\u00d73=1263\", \"526\u00d77=3682\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the text replacements described by the diff: the header date and\n# the 25 three-digit-by-one-digit multiplication equations in the table.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2025-03-15 Saturday\", \"2025-03-16 Sunday\"),\n  @(\"229\u00d72=458\", \"457\u00d74=1828\"),\n  @(\"981\u00d77=6867\", \"988\u00d75=4940\"),\n  @(\"730\u00d76=4380\", \"393\u00d77=2751\"),\n  @(\"588\u00d79=5292\", \"157\u00d74=628\"),\n  @(\"324\u00d72=648\", \"357\u00d72=714\"),\n  @(\"815\u00d75=4075\", \"143\u00d78=1144\"),\n  @(\"920\u00d72=1840\", \"916\u00d75=4580\"),\n  @(\"867\u00d72=1734\", \"461\u00d74=1844\"),\n  @(\"909\u00d72=1818\", \"458\u00d75=2290\"),\n  @(\"765\u00d77=5355\", \"616\u00d75=3080\"),\n  @(\"654\u00d72=1308\", \"859\u00d76=5154\"),\n  @(\"515\u00d73=1545\", \"329\u00d72=658\"),\n  @(\"238\u00d77=1666\", \"425\u00d72=850\"),\n  @(\"450\u00d72=900\", \"216\u00d78=1728\"),\n  @(\"805\u00d73=2415\", \"245\u00d73=735\"),\n  @(\"277\u00d75=1385\", \"276\u00d77=1932\"),\n  @(\"291\u00d73=873\", \"344\u00d76=2064\"),\n  @(\"651\u00d73=1953\", \"136\u00d77=952\"),\n  @(\"116\u00d77=812\", \"516\u00d78=4128\"),\n  @(\"472\u00d72=944\", \"337\u00d73=1011\"),\n  @(\"628\u00d76=3768\", \"188\u00d75=940\"),\n  @(\"496\u00d79=4464\", \"136\u00d74=544\"),\n  @(\"227\u00d77=1589\", \"972\u00d78=7776\"),\n  @(\"582\u00d76=3492\", \"920\u00d77=6440\"),\n  @(\"421\u00d73=1263\", \"526\u00d77=3682\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
